$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 7 (line "6" in the logic table) currently has "4" and "5" in the
# "First Segment" / "Second Segment" columns; clear that text, leaving
# the paragraphs empty (no run), matching a revert of the earlier edit.
$row6 = $t.Rows.Item(7)

$cell4 = $row6.Cells.Item(4)
$rng4 = $cell4.Range
$d.Range($rng4.Start, $rng4.Start + 1).Delete()

$cell5 = $row6.Cells.Item(5)
$rng5 = $cell5.Range
$d.Range($rng5.Start, $rng5.Start + 1).Delete()

# Row 8 (line "7": ¬p / ¬i / 3-6) is entirely removed.
$t.Rows.Item(8).Delete()
